$wb = $excel.ActiveWorkbook

# Sheet "Imm.0.8.3": append "*" to gene names in A2:A11
$ws1 = $wb.Worksheets.Item("Imm.0.8.3")
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + "*"
}

# Sheet "Imm.1.2.14": append "*" to gene names in A3:A13
$ws2 = $wb.Worksheets.Item("Imm.1.2.14")
for ($r = 3; $r -le 13; $r++) {
    $cell = $ws2.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + "*"
}
